$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.848583221435547
$ws.Range("B1").Value = 5.535118579864502
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.818385362625122
$ws.Range("E1").Value = 1.964761137962341
